$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N4").Value = 100
$ws.Range("N5").Value = 100
$ws.Range("N6").Value = 100
$ws.Range("N7").Value = 100
$ws.Range("N8").Value = 100

$ws.Range("P9").Select()
